$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C13").Value = "Admin_Report()"
$ws.Range("F16").Select()
